$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CurrProg_Summary")

$ws.Range("B2").Value = 3736.72709330646
$ws.Range("C2").Value = 2926033

$ws.Range("B3").Value = 5223.53644314869
$ws.Range("C3").Value = 3179840

$ws.Range("B4").Value = 43068.5287981859
$ws.Range("C4").Value = 6235918

$ws.Range("B5").Value = 936.86658286503
$ws.Range("C5").Value = 1248880
